# Update the workbook: append two new rows (76 and 77) of data to the
# "Optical_Power" worksheet, matching the records added in the source diff.
#
# Values that look numeric/date-like (Caso, F. De Reclamo, Comuna, OT,
# Attachments) are entered with a leading apostrophe so Excel keeps them
# as plain text, matching the rest of the sheet where these columns are
# stored as text rather than numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76
$ws.Range("A76").Value = "'6357"
$ws.Range("B76").Value = "'7/7/2025"
$ws.Range("C76").Value = "BACACAY 3088"
$ws.Range("D76").Value = "'7"
$ws.Range("E76").Value = "'808036196"
$ws.Range("F76").Value = "Optical Power"
$ws.Range("G76").Value = "Pendiente"
$ws.Range("H76").Value = "Reparar rienda"
$ws.Range("I76").Value = "'1"
$ws.Range("J76").Value = "Tensor"
$ws.Range("K76").Value = "Sin equipos"
$ws.Range("L76").Value = "Terminal"

# Row 77
$ws.Range("A77").Value = "'-502"
$ws.Range("B77").Value = "'7/7/2025"
$ws.Range("C77").Value = "Tagle 2562"
$ws.Range("D77").Value = "'14"
$ws.Range("E77").Value = "'808036198"
$ws.Range("F77").Value = "Optical Power"
$ws.Range("G77").Value = "Pendiente"
$ws.Range("H77").Value = "Colocar columna para pedir traspaso nodo teco"
$ws.Range("I77").Value = "'1"
$ws.Range("J77").Value = "Cambio"
$ws.Range("K77").Value = "Nodo Teco"
$ws.Range("L77").Value = "Pasante"
